$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Fill in the NRF24l01 pinout row for the Uno (row 5, columns O:U), which
#    was previously left blank (only Nano and Mega rows were populated).
# ---------------------------------------------------------------------------
$ws.Range("O5").Value2 = "3V3"
$ws.Range("P5").Value2 = "-"
$ws.Range("Q5").Value2 = "D8"
$ws.Range("R5").Value2 = "D7"
$ws.Range("S5").Value2 = "D11 - MOSI"
$ws.Range("T5").Value2 = "D13 - SCK"
$ws.Range("U5").Value2 = "D12 - MISO"

$srcRow4 = $ws.Range("O4:U4")
$dstRow5 = $ws.Range("O5:U5")
$dstRow5.HorizontalAlignment = $srcRow4.HorizontalAlignment
$dstRow5.VerticalAlignment = $srcRow4.VerticalAlignment

# ---------------------------------------------------------------------------
# 2. Add a new "Servo" accessory section in columns Z:AB with the Uno pinout.
# ---------------------------------------------------------------------------
$ws.Range("Z2:AB2").Merge()
$ws.Range("Z2").Value2 = "Servo"
$ws.Range("Z2:AB2").HorizontalAlignment = $ws.Range("V2").HorizontalAlignment
$ws.Range("Z2:AB2").VerticalAlignment = $ws.Range("V2").VerticalAlignment

$ws.Range("Z3").Value2 = "+"
$ws.Range("AA3").Value2 = "-"
$ws.Range("AB3").Value2 = "Data"
$ws.Range("Z3:AB3").HorizontalAlignment = $ws.Range("V3").HorizontalAlignment
$ws.Range("Z3:AB3").VerticalAlignment = $ws.Range("V3").VerticalAlignment

$ws.Range("Z5").Value2 = "5V"
$ws.Range("AA5").Value2 = "GND"
$ws.Range("AB5").Value2 = "D9"
$ws.Range("Z5:AB5").HorizontalAlignment = $ws.Range("V4").HorizontalAlignment
$ws.Range("Z5:AB5").VerticalAlignment = $ws.Range("V4").VerticalAlignment

# ---------------------------------------------------------------------------
# 3. Update Mega's NRF24l01 MOSI/SCK/MISO pin labels to match the datasheet
#    naming convention (prefix pin numbers with "D").
# ---------------------------------------------------------------------------
$ws.Range("S6").Value2 = "D51 - MOSI"
$ws.Range("T6").Value2 = "D52 - SCK"
$ws.Range("U6").Value2 = "D50 - MISO"

# ---------------------------------------------------------------------------
# 4. Update the selection / view state to reflect where the edit happened.
# ---------------------------------------------------------------------------
$ws.Range("X13").Select()
